# Chapter13-HW-3.xlsx edit script
#
# 1. Rename the sheet "2-hiddenLayers" -> "NN"
#    (Excel automatically rewrites every localSheetId-scoped defined name
#    that referenced '2-hiddenLayers' to point at NN instead.)
# 2. Update the instructions text box:
#      "Neuron Network procedure" -> "Neural Network (NN) procedure"
#      "Note, using one or two hidden layers is acceptable. But it would be
#       interesting to test the result of using 3 hidden layers"
#        -> "Note, using one or two hidden layers only."
# 3. Move the active-cell selection from I23 to I17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet -------------------------------------------------
$ws.Name = "NN"

# --- 2. Update the floating text box instructions ----------------------------
$shp = $ws.Shapes.Item(1)

$crlf = [char]13
$newText = "Given the training dataset below, please use Neural Network (NN) procedure with 4 neurons to construst a NN model and apply the model to the scoring dataset (below the training dataset)." + $crlf + $crlf + "Note, using one or two hidden layers only."

$shp.TextFrame.Characters().Text = $newText

# --- 3. Update the selected cell ---------------------------------------------
$ws.Range("I17").Select() | Out-Null
